# Add the 2020-03-31 data row (row 43) to the COVID-19 infection dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43
$ws.Cells.Item($row, 1).Value = 43920
$ws.Cells.Item($row, 2).Value = 6663
$ws.Cells.Item($row, 3).Value = 1414
$ws.Cells.Item($row, 4).Value = 75
$ws.Cells.Item($row, 5).Value = 122
$ws.Cells.Item($row, 6).Value = 5249
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0

# Column A carries a yyyy-mm-dd date format; make sure the new date cell matches.
$ws.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd"

# Scroll the view down to show the newly added rows, matching the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("H42").Select()
